$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.591.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.418.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.459.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0986"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.861.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.478.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.493.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "316.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.559.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.385"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.154"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0739"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.992"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.825"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.30%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "261.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.575"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0920"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0499"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0215"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.96%  "
